$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B is the "Rule" name/label column (see the cell comment on B4:
# "Rule column is used to name particular rule rows for documentation and
# tracing purposes"). The rule previously labeled "R40" in B11 is renamed
# to "1". A leading apostrophe forces Excel to store the literal as text
# (a shared string), matching the original label's text type, rather than
# silently re-interpreting the digit as a number.
$ws.Range("B11").Value = "'1"
